$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

foreach ($ws in @($ws1, $ws4)) {
    $ws.Range("F2").Value = 455

    $ws.Range("C3").Value = "丽水·CCAC动漫游戏嘉年华"
    $ws.Range("D3").Value = "南环西路109号 九城宴会中心"
    $ws.Range("E3").Value = "2024.07.20 09:00-07.20 16:00"
    $ws.Range("F3").Value = 17
    $ws.Range("G3").Value = 29.9
    $ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=86306"
    $ws.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202405/4TrBjBlV1716551375116.png"

    $ws.Range("C4").Value = "龙泉·ACG动湿游戏博览会"
    $ws.Range("D4").Value = "南秦路1号望瓯·陶溪川直走200米左手边(7号楼) 望瓯陶溪川活动中心"
    $ws.Range("E4").Value = "2024.07.20 10:00-07.21 18:00"
    $ws.Range("F4").Value = 1
    $ws.Range("G4").Value = 55
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=86671"
    $ws.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202405/sg6nrCrJ1717142810026.png"
}
